$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 4 (sCs target), shifting old row4 -> row5
$ws.Rows.Item(4).Insert()

# --- Row 2 (sCs -> ECs) updates ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.180598
$ws.Range("H2").Value = 18.541794
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.239942333333333
$ws.Range("N2").Value = 3.719827
$ws.Range("O2").Value = 0.02873264962564676
$ws.Range("P2").Value = 0.02873264962564676
$ws.Range("Q2").Value = 7.663585105515333
$ws.Range("R2").Value = 68.97226594963799
$ws.Range("S2").Value = 0.02873264962564676
$ws.Range("T2").Value = 0.02873264962564676

# --- Row 3 (sCs -> FAPs) updates ---
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.180598
$ws.Range("H3").Value = 18.541794
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.025359
$ws.Range("N3").Value = 81.076077
$ws.Range("O3").Value = 0.6262470038157576
$ws.Range("P3").Value = 0.6262470038157575
$ws.Range("Q3").Value = 167.032879784682
$ws.Range("R3").Value = 1503.295918062138
$ws.Range("S3").Value = 0.6262470038157576
$ws.Range("T3").Value = 0.6262470038157575

# --- Row 4 (new row, sCs -> M2) ---
$ws.Range("A4").Value2 = "sCs"
$ws.Range("B4").Value2 = "Edn3"
$ws.Range("C4").Value2 = "Ednra"
$ws.Range("D4").Value2 = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.180598
$ws.Range("H4").Value = 18.541794
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.050758
$ws.Range("N4").Value = 0.152274
$ws.Range("O4").Value = 0.001176193271648315
$ws.Range("P4").Value = 0.001176193271648314
$ws.Range("Q4").Value = 0.313714793284
$ws.Range("R4").Value = 2.823433139556
$ws.Range("S4").Value = 0.001176193271648315
$ws.Range("T4").Value = 0.001176193271648314

# --- Row 5 (previously row 4, sCs -> sCs) updates ---
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.180598
$ws.Range("H5").Value = 18.541794
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.838413
$ws.Range("N5").Value = 44.51523899999999
$ws.Range("O5").Value = 0.3438441532869475
$ws.Range("P5").Value = 0.3438441532869475
$ws.Range("Q5").Value = 91.71026571097399
$ws.Range("R5").Value = 825.3923913987659
$ws.Range("S5").Value = 0.3438441532869475
$ws.Range("T5").Value = 0.3438441532869475
